# Adds a new "End Turn" localization row (GAME / game.endturn / <translation>)
# to each of the four language sheets, then restores the UI navigation state
# (active sheet + selections) to match the post-edit workbook.

$wb = $excel.ActiveWorkbook

# sheet name -> (row number to append, translated "End Turn" text)
$rows = @{
    "English"    = @{ Row = 57; Text = "End Turn" }
    "Svenska"    = @{ Row = 45; Text = "Avsluta Tur" }
    "Ελληνικά"   = @{ Row = 45; Text = "Τέλος στροφής" }
    "Nederlands" = @{ Row = 45; Text = "Einde Beurt" }
}

foreach ($ws in $wb.Worksheets) {
    $info = $rows[$ws.Name]
    if ($info -eq $null) { continue }

    $r = $info.Row

    $catCell = $ws.Cells.Item($r, 1)
    $catCell.Value = "GAME"
    $catCell.Font.Bold = $true

    $keyCell = $ws.Cells.Item($r, 2)
    $keyCell.Value = "game.endturn"

    $txtCell = $ws.Cells.Item($r, 3)
    $txtCell.Value = $info.Text
    # Column C on the other data rows inherits a column-level style; the new
    # row's text cell instead uses the workbook's plain default formatting,
    # so pin the font explicitly back to the default Calibri 12 rather than
    # letting it inherit the neighbouring rows' Arial 10 style.
    $txtCell.Font.Name = "Calibri"
    $txtCell.Font.Size = 12
}

# Make "Svenska" the active sheet/tab (was "English").
$svenska = $wb.Worksheets.Item("Svenska")
$svenska.Activate()
$svenska.Range("C46").Select() | Out-Null

# Update the other sheets' remembered selections to their new last rows.
$english = $wb.Worksheets.Item("English")
$english.Range("A57").Select() | Out-Null

$greek = $wb.Worksheets.Item("Ελληνικά")
$greek.Range("C45").Select() | Out-Null

# Re-activate Svenska so it ends up as the tab shown/selected on open.
$svenska.Activate()
$svenska.Range("C46").Select() | Out-Null
